$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Youssef")

# Fill in the new "Delivery Guy" user-story rows (3-24) in the
# Module / Depends-on / Requirement ID / As a / I want to / So that / Platform table.
$ws.Cells.Item(3,3).Value = 1
$ws.Cells.Item(3,4).Value = "Delivery Guy"
$ws.Cells.Item(3,5).Value = "be able to sign in as a delivery guy"
$ws.Cells.Item(3,6).Value = "I can view my delivery portal"
$ws.Cells.Item(3,7).Value = "Both"
$ws.Cells.Item(4,1).Value = "Delivery"
$ws.Cells.Item(4,2).Value = "-"
$ws.Cells.Item(4,3).Value = 2
$ws.Cells.Item(4,4).Value = "Delivery Guy"
$ws.Cells.Item(4,5).Value = "be able to view orders ready for delivery"
$ws.Cells.Item(4,6).Value = "I can choose from them"
$ws.Cells.Item(4,7).Value = "Both"
$ws.Cells.Item(5,1).Value = "Delivery"
$ws.Cells.Item(5,2).Value = "-"
$ws.Cells.Item(5,3).Value = 3
$ws.Cells.Item(5,4).Value = "Delivery Guy"
$ws.Cells.Item(5,5).Value = "be able to select an order"
$ws.Cells.Item(5,6).Value = "no other delivery guy can select it"
$ws.Cells.Item(5,7).Value = "Both"
$ws.Cells.Item(6,1).Value = "Delivery"
$ws.Cells.Item(6,2).Value = "-"
$ws.Cells.Item(6,3).Value = 4
$ws.Cells.Item(6,4).Value = "Delivery Guy"
$ws.Cells.Item(6,5).Value = "be able to see donor and receiver details"
$ws.Cells.Item(6,6).Value = "I can reach them"
$ws.Cells.Item(6,7).Value = "Both"
$ws.Cells.Item(7,1).Value = "Delivery"
$ws.Cells.Item(7,2).Value = "-"
$ws.Cells.Item(7,3).Value = 5
$ws.Cells.Item(7,4).Value = "Delivery Guy"
$ws.Cells.Item(7,5).Value = "be able to open donor's and receiver's locations"
$ws.Cells.Item(7,6).Value = "I can reach them"
$ws.Cells.Item(7,7).Value = "Both"
$ws.Cells.Item(8,1).Value = "Delivery"
$ws.Cells.Item(8,2).Value = "-"
$ws.Cells.Item(8,3).Value = 6
$ws.Cells.Item(8,4).Value = "Delivery Guy"
$ws.Cells.Item(8,5).Value = "be able to know the dimensions and weight of the object"
$ws.Cells.Item(8,6).Value = "I can know if i can transport it"
$ws.Cells.Item(8,7).Value = "Both"
$ws.Cells.Item(9,1).Value = "Delivery"
$ws.Cells.Item(9,2).Value = "-"
$ws.Cells.Item(9,3).Value = 7
$ws.Cells.Item(9,4).Value = "Delivery Guy"
$ws.Cells.Item(9,5).Value = "be able to know if the package has fragile items"
$ws.Cells.Item(9,6).Value = "I can be careful when transporting it"
$ws.Cells.Item(9,7).Value = "Both"
$ws.Cells.Item(10,1).Value = "Delivery"
$ws.Cells.Item(10,2).Value = "-"
$ws.Cells.Item(10,3).Value = 8
$ws.Cells.Item(10,4).Value = "Delivery Guy"
$ws.Cells.Item(10,5).Value = "be able to know if the package needs to be cooled"
$ws.Cells.Item(10,6).Value = "I can careful when transporting it"
$ws.Cells.Item(10,7).Value = "Both"
$ws.Cells.Item(11,1).Value = "Delivery"
$ws.Cells.Item(11,2).Value = "-"
$ws.Cells.Item(11,3).Value = 9
$ws.Cells.Item(11,4).Value = "Delivery Guy"
$ws.Cells.Item(11,5).Value = "be able to view collecting donation time"
$ws.Cells.Item(11,6).Value = "I know when to go to the donor"
$ws.Cells.Item(11,7).Value = "Both"
$ws.Cells.Item(12,1).Value = "Delivery"
$ws.Cells.Item(12,2).Value = "-"
$ws.Cells.Item(12,3).Value = 10
$ws.Cells.Item(12,4).Value = "Delivery Guy"
$ws.Cells.Item(12,5).Value = "be able to view receiving time of donation"
$ws.Cells.Item(12,6).Value = "I know when to go to the receiver"
$ws.Cells.Item(12,7).Value = "Both"
$ws.Cells.Item(13,1).Value = "Delivery"
$ws.Cells.Item(13,2).Value = "-"
$ws.Cells.Item(13,3).Value = 11
$ws.Cells.Item(13,4).Value = "Delivery Guy"
$ws.Cells.Item(13,5).Value = "be able to receive money through the application"
$ws.Cells.Item(13,6).Value = "I receive wages"
$ws.Cells.Item(13,7).Value = "Both"
$ws.Cells.Item(14,1).Value = "Delivery"
$ws.Cells.Item(14,2).Value = "-"
$ws.Cells.Item(14,3).Value = 12
$ws.Cells.Item(14,4).Value = "Delivery Guy"
$ws.Cells.Item(14,5).Value = "be able to tick that I received wages through the application"
$ws.Cells.Item(14,6).Value = "wages delivery is monitored"
$ws.Cells.Item(14,7).Value = "Both"
$ws.Cells.Item(15,1).Value = "Delivery"
$ws.Cells.Item(15,2).Value = "-"
$ws.Cells.Item(15,3).Value = 13
$ws.Cells.Item(15,4).Value = "Delivery Guy"
$ws.Cells.Item(15,5).Value = "be able to tick that I delivered the donation"
$ws.Cells.Item(15,6).Value = "packages delivery is monitored"
$ws.Cells.Item(15,7).Value = "Both"
$ws.Cells.Item(16,1).Value = "Delivery"
$ws.Cells.Item(16,2).Value = "-"
$ws.Cells.Item(16,3).Value = 14
$ws.Cells.Item(16,4).Value = "Delivery Guy"
$ws.Cells.Item(16,5).Value = "have access to a customer feedback system within the app"
$ws.Cells.Item(16,6).Value = "I can receive feedback on my performance and make improvements if necessary"
$ws.Cells.Item(16,7).Value = "Both"
$ws.Cells.Item(17,1).Value = "Delivery"
$ws.Cells.Item(17,2).Value = "-"
$ws.Cells.Item(17,3).Value = 15
$ws.Cells.Item(17,4).Value = "Delivery Guy"
$ws.Cells.Item(17,5).Value = "receive timely notifications and alerts from the delivery app regarding any changes to delivery schedules or new delivery assignments"
$ws.Cells.Item(17,6).Value = "I can adapt my plans accordingly"
$ws.Cells.Item(17,7).Value = "Both"
$ws.Cells.Item(18,1).Value = "Delivery"
$ws.Cells.Item(18,2).Value = "-"
$ws.Cells.Item(18,3).Value = 16
$ws.Cells.Item(18,4).Value = "Delivery Guy"
$ws.Cells.Item(18,5).Value = "have access to a rewards or incentive program within the app"
$ws.Cells.Item(18,6).Value = "I can earn rewards or bonuses based on my performance"
$ws.Cells.Item(18,7).Value = "Both"
$ws.Cells.Item(19,1).Value = "Delivery"
$ws.Cells.Item(19,2).Value = "-"
$ws.Cells.Item(19,3).Value = 17
$ws.Cells.Item(19,4).Value = "Delivery Guy"
$ws.Cells.Item(19,5).Value = "have access to a history of past deliveries within the app"
$ws.Cells.Item(19,6).Value = "I can easily reference previous delivery details and track my performance over time"
$ws.Cells.Item(19,7).Value = "Both"
$ws.Cells.Item(20,1).Value = "Delivery"
$ws.Cells.Item(20,2).Value = "-"
$ws.Cells.Item(20,3).Value = 18
$ws.Cells.Item(20,4).Value = "Delivery Guy"
$ws.Cells.Item(20,5).Value = "have access to a feedback mechanism within the app"
$ws.Cells.Item(20,6).Value = "I can provide feedback on the overall delivery process, app functionality, and any suggestions for improvement"
$ws.Cells.Item(20,7).Value = "Both"
$ws.Cells.Item(21,1).Value = "Delivery"
$ws.Cells.Item(21,2).Value = "-"
$ws.Cells.Item(21,3).Value = 19
$ws.Cells.Item(21,4).Value = "Delivery Guy"
$ws.Cells.Item(21,5).Value = "have access to a feature within the app that allows me to easily report any vehicle maintenance issues or concerns"
$ws.Cells.Item(21,6).Value = "they can be addressed promptly to ensure the safety and reliability of the delivery vehicle"
$ws.Cells.Item(21,7).Value = "Both"
$ws.Cells.Item(22,1).Value = "Delivery"
$ws.Cells.Item(22,2).Value = "-"
$ws.Cells.Item(22,3).Value = 20
$ws.Cells.Item(22,4).Value = "Delivery Guy"
$ws.Cells.Item(22,5).Value = "have access to a feature within the app that allows me to set reminders for pickup times and delivery deadlines"
$ws.Cells.Item(22,6).Value = "I can stay organized and on schedule"
$ws.Cells.Item(22,7).Value = "Both"
$ws.Cells.Item(23,1).Value = "Delivery"
$ws.Cells.Item(23,2).Value = "-"
$ws.Cells.Item(23,3).Value = 21
$ws.Cells.Item(23,4).Value = "Delivery Guy"
$ws.Cells.Item(23,5).Value = "have access to a feature within the app that allows me to track my earnings and incentives in real-time"
$ws.Cells.Item(23,6).Value = "I can ensure accurate and transparent compensation for my work"
$ws.Cells.Item(23,7).Value = "Both"
$ws.Cells.Item(24,1).Value = "Delivery"
$ws.Cells.Item(24,2).Value = "-"
$ws.Cells.Item(24,3).Value = 22
$ws.Cells.Item(24,4).Value = "Delivery Guy"
$ws.Cells.Item(24,5).Value = "have access to a feature within the app that allows me to report and document any instances of package theft, damage, or tampering during delivery"
$ws.Cells.Item(24,6).Value = "investigations are facilitated and accountability is ensured"
$ws.Cells.Item(24,7).Value = "Both"

# Widen columns E and F so the longer story/rationale text added above is readable.
$ws.Columns.Item(5).ColumnWidth = 119.16666666666667
$ws.Columns.Item(6).ColumnWidth = 87.75

# Leave the selection where the author ended up after entering the data.
$ws.Range("D28").Select()
